$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O11").Value = 15
$ws.Range("T11").Value = 21

$ws.Range("A12").Value = "lorena.caixeta@mrv.com.br"
$ws.Range("F12").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004"
$ws.Range("I12").Value = "2025-05-20 11:13:31"
$ws.Range("J12").Value = "Painel do Portifólio - Planejamento da Produção - PLNESROBR004: ok."
$ws.Range("K12").Value = "Planilha de balanceamento de carga horária"
$ws.Range("M12").Value = "Gestão de Projetos"
$ws.Range("N12").Value = 1
$ws.Range("O12").NumberFormat = "@"
$ws.Range("O12").Value = "12"
